$d = $word.ActiveDocument

# Helper: find a paragraph whose visible text (paragraph mark stripped)
# equals the given target string exactly.
function Find-ParagraphByExactText($doc, $targetText) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $p
        }
    }
    return $null
}

# Helper: append one or more new runs (each its own <w:r>) to the end of a
# paragraph's text (before its paragraph mark), using Range.InsertXML with a
# Flat-OPC package so that every supplied run stays a distinct <w:r> element
# instead of being merged into the paragraph's existing run.
function Add-RunsToParagraphEnd($paragraph, $runsXml) {
    $r = $paragraph.Range
    $r.MoveEnd(1, -1) # exclude the trailing paragraph mark
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1) "Ổ cứng:" paragraph -> append " " + "(cho option chọn )"
$pOCung = Find-ParagraphByExactText $d "Ổ cứng:"
if ($pOCung -ne $null) {
    $runs = '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(cho option chọn )</w:t></w:r>'
    Add-RunsToParagraphEnd $pOCung $runs
}

# 2) "Ram" paragraph -> append ": " + "(cho option chọn )"
$pRam = Find-ParagraphByExactText $d "Ram"
if ($pRam -ne $null) {
    $runs = '<w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>(cho option chọn )</w:t></w:r>'
    Add-RunsToParagraphEnd $pRam $runs
}
